# Rename Sheet2 -> DPdata
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "DPdata"

# Populate DPdata (formerly Sheet2) with data
$ws2.Range("A1").Value = "jisha"
$ws2.Range("B1").Value = "admin"
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "admin"
$ws2.Range("A3").Value = "ghhg"
$ws2.Range("B3").Value = "jhj"

# Update selection on Sheet1: selected range A2:B3 (active cell B3)
$ws1.Range("A2:B3").Select()

# Make DPdata the active sheet, with selection at B4
$ws2.Activate()
$ws2.Range("B4").Select()

$wb.Save()
